$d = $word.ActiveDocument

# --- 1) Fix the two VML "w:pict" runs: add a w14:anchorId and renumber the
#        legacy VML shape id (Word also bumps o:title-less "Capture N" shape
#        ids when it re-saves pasted/linked images). We locate the owning
#        paragraphs by sniffing their OOXML for the old shape id, then do a
#        surgical WordOpenXML round trip (read -> string edit -> InsertXML)
#        so only that paragraph's markup is touched.
function Fix-PictParagraph($doc, $oldShapeId, $newShapeId, $anchorId) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        $openXml = $para.Range.WordOpenXML
        if ($openXml.Contains($oldShapeId)) {
            $updated = $openXml.Replace("<w:pict>", "<w:pict w14:anchorId=`"$anchorId`">")
            $updated = $updated.Replace($oldShapeId, $newShapeId)
            # Drop the paraId/textId stamp the read-back snapshot carries so
            # InsertXML doesn't bake a fresh (Word-2010-era) paragraph id
            # into a document that never used them.
            $updated = [regex]::Replace($updated, ' w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+"', '')
            $para.Range.InsertXML($updated)
            return $true
        }
    }
    return $false
}

Fix-PictParagraph $d "_x0000_i1029" "_x0000_i1025" "633AD93F" | Out-Null
Fix-PictParagraph $d "_x0000_i1030" "_x0000_i1026" "1C02C124" | Out-Null
